$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "97.045.66"
$ws.Range("E2").Value = "  +1.78%  "

# Row 3
$ws.Range("D3").Value = "3.567.56"
$ws.Range("E3").Value = "  -0.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'241.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "

# Row 6
$ws.Range("D6").Value = "'655.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.37%  "

# Row 7
$ws.Range("E7").Value = "  +15.69%  "

# Row 8
$ws.Range("E8").Value = "  +6.51%  "

# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").Value = "'1.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.79%  "

# Row 11
$ws.Range("D11").Value = "3.563.40"
$ws.Range("E11").Value = "  -0.47%  "

# Row 12
$ws.Range("D12").Value = "'44.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.46%  "

# Row 13
$ws.Range("E13").Value = "  +0.30%  "

# Row 14
$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("D15").Value = "4.234.76"
$ws.Range("E15").Value = "  -0.26%  "

# Row 16
$ws.Range("D16").Value = "97.082.95"
$ws.Range("E16").Value = "  +1.93%  "

# Row 17
$ws.Range("D17").Value = "'0.0000259"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.20%  "

# Row 18
$ws.Range("D18").Value = "'8.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +10.93%  "

# Row 19
$ws.Range("D19").Value = "3.568.73"
$ws.Range("E19").Value = "  -0.29%  "

# Row 20
$ws.Range("D20").Value = "'12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.93%  "

# Row 21
$ws.Range("D21").Value = "'17.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

# Row 22
$ws.Range("D22").Value = "'0.524"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.57%  "

# Row 23
$ws.Range("D23").Value = "'3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.12%  "

# Row 24
$ws.Range("D24").Value = "'512.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "

# Row 25
$ws.Range("D25").Value = "'0.0000204"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.40%  "

# Row 26
$ws.Range("D26").Value = "'6.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "

# Row 27
$ws.Range("D27").Value = "'101.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.95%  "

# Row 28
$ws.Range("D28").Value = "'12.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.12%  "

# Row 29
$ws.Range("D29").Value = "3.763.06"
$ws.Range("E29").Value = "  -0.28%  "

# Row 30
$ws.Range("D30").Value = "'0.162"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.78%  "

# Row 31
$ws.Range("D31").Value = "'2.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

# Row 32
$ws.Range("D32").Value = "'11.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.49%  "

# Row 33
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.47%  "

# Row 34
$ws.Range("E34").Value = "  +3.33%  "

# Row 35
$ws.Range("D35").Value = "'0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.41%  "

# Row 36
$ws.Range("D36").Value = "'31.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("D37").Value = "'8.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "

# Row 38
$ws.Range("D38").Value = "'614.79"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.61%  "

# Row 39
$ws.Range("D39").Value = "'0.564"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "

# Row 40
$ws.Range("D40").Value = "'1.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.154"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.36%  "

# Row 42
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").Value = "'1.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.65%  "

# Row 43
$ws.Range("E43").Value = "  -0.09%  "

# Row 44
$ws.Range("D44").Value = "'0.917"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.12%  "

# Row 45
$ws.Range("D45").Value = "'5.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.06%  "

# Row 46
$ws.Range("D46").Value = "'0.0436"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.02%  "

# Row 47
$ws.Range("D47").Value = "'2.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48
$ws.Range("D48").Value = "'23.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.83%  "

# Row 49
$ws.Range("D49").Value = "'0.401"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +27.75%  "

# Row 50
$ws.Range("D50").Value = "'8.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.37%  "

# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'32.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.14%  "
